$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values (swap with old row 5 content)
$ws.Range("D2").Value = 44181
$ws.Range("M2").Value = 65
$ws.Range("N2").Value = 3600
$ws.Range("O2").Value = 3800
$ws.Range("P2").Value = 3692
$ws.Range("R2").Value = "Provincia de Diguillín"
$ws.Range("S2").Value = 1846

# Row 3 values (swap with old row 6 content)
$ws.Range("D3").Value = 44181
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 1800
$ws.Range("O3").Value = 2000
$ws.Range("P3").Value = 1875
$ws.Range("S3").Value = 1875

# Row 5 values (swap with old row 2 content)
$ws.Range("D5").Value = 44187
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 2800
$ws.Range("O5").Value = 3000
$ws.Range("P5").Value = 2900
$ws.Range("R5").Value = "Provincia de Linares"
$ws.Range("S5").Value = 1450

# Row 6 values (swap with old row 3 content)
$ws.Range("D6").Value = 44187
$ws.Range("M6").Value = 65
$ws.Range("N6").Value = 1400
$ws.Range("O6").Value = 1500
$ws.Range("P6").Value = 1446
$ws.Range("S6").Value = 1446
